$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.805.18"
$ws.Range("E2").Value = "  -7.35%  "

# Row 3
$ws.Range("D3").Value = "1.701.99"
$ws.Range("E3").Value = "  -6.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.70"
$ws.Range("E5").Value = "  -5.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5141"
$ws.Range("E6").Value = "  -13.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  -6.20%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "22.04"
$ws.Range("E9").Value = "  -3.65%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06193"
$ws.Range("E10").Value = "  -8.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07337"
$ws.Range("E11").Value = "  -2.12%  "

# Row 12
$ws.Range("D12").Value = "1.709.05"
$ws.Range("E12").Value = "  -6.66%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.478"
$ws.Range("E13").Value = "  -4.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5801"
$ws.Range("E14").Value = "  -6.92%  "

# Row 15
$ws.Range("D15").Value = "1.934.44"
$ws.Range("E15").Value = "  -6.11%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008182"
$ws.Range("E16").Value = "  -11.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.39"

# Row 18
$ws.Range("D18").Value = "26.884.77"
$ws.Range("E18").Value = "  -6.04%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.002"
$ws.Range("E19").Value = "  -7.95%  "

# Row 20
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.65"
$ws.Range("E21").Value = "  -6.37%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "186.12"
$ws.Range("E22").Value = "  -10.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.262"
$ws.Range("E23").Value = "  -7.42%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.18"
$ws.Range("E25").Value = "  -7.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.555"
$ws.Range("E26").Value = "  -2.98%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1153"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.21"
$ws.Range("E28").Value = "  -6.69%  "

# Row 29
$ws.Range("E29").Value = "  -4.55%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05863"
$ws.Range("E30").Value = "  -7.28%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.346"
$ws.Range("E31").Value = "  -5.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.460"
$ws.Range("E32").Value = "  -7.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.437"
$ws.Range("E33").Value = "  -6.62%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.646"
$ws.Range("E34").Value = "  -2.51%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9865"
$ws.Range("E35").Value = "  -5.86%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.412"
$ws.Range("E36").Value = "  -4.32%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6013"
$ws.Range("E37").Value = "  -4.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.675"
$ws.Range("E38").Value = "  -2.08%  "

# Row 39
$ws.Range("D39").Value = "1.098.06"
$ws.Range("E39").Value = "  -2.97%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01598"
$ws.Range("E40").Value = "  -5.35%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.930"
$ws.Range("E41").Value = "  -7.62%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8634"
$ws.Range("E42").Value = "  -0.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.53"
$ws.Range("E44").Value = "  -2.49%  "

# Row 45
$ws.Range("D45").Value = "1.848.79"
$ws.Range("E45").Value = "  -6.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.11"
$ws.Range("E46").Value = "  -7.06%  "

# Row 47
$ws.Range("E47").Value = "  -6.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9998"
$ws.Range("E48").Value = "  -0.29%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4353"
$ws.Range("E49").Value = "  -3.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05243"
$ws.Range("E50").Value = "  -3.98%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.965"
$ws.Range("E51").Value = "  -3.75%  "
